$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C3").Value = "Is this the real sheet"
